# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right before the existing "2021-Q4"
#   sheet (so the tab order becomes 总计, 2022-Q3, 2021-Q4, 2021-Q2).
# - Populate it with the 2022-Q3 fund-holding table.
# - Insert a matching summary row into the "总计" sheet, pushing the
#   previously-existing rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "2022-Q3" worksheet, positioned before "2021-Q4".
# ---------------------------------------------------------------------------
$q4Before = $wb.Worksheets.Item(2)          # "2021-Q4", currently 2nd tab
$wb.Worksheets.Add($q4Before) | Out-Null    # Add() with a target inserts *before* it

$newSheet = $wb.Worksheets.Item(2)          # the freshly inserted sheet
$newSheet.Name = "2022-Q3"

$oldQ4 = $wb.Worksheets.Item(3)             # "2021-Q4" shifted to the 3rd tab

# ---------------------------------------------------------------------------
# Copy the formatting (borders/bold/alignment) used on the "2021-Q4" sheet's
# header row and index column onto the same cells of the new sheet, so the
# new sheet looks consistent with its siblings.
# ---------------------------------------------------------------------------
$oldQ4.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$oldQ4.Range("A2").Copy($newSheet.Range("A2"))
$oldQ4.Range("A2").Copy($newSheet.Range("A3"))
$oldQ4.Range("A2").Copy($newSheet.Range("A4"))
$oldQ4.Range("A2").Copy($newSheet.Range("A5"))

# ---------------------------------------------------------------------------
# Header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# Data rows. Columns B, D, E, F, G hold numeric-looking text (fund codes /
# percentages kept as strings, e.g. so leading zeros survive), so force a
# text number-format before writing them. Column A (index) and H (rank)
# are genuine numbers.
# ---------------------------------------------------------------------------
# NOTE: comma-separated multi-area Range() strings are not reliably honoured
# by this COM host, so each contiguous block is addressed separately.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "213003"
$newSheet.Range("C2").Value = "宝盈策略增长混合"
$newSheet.Range("D2").Value = "10.55"
$newSheet.Range("E2").Value = "90.74"
$newSheet.Range("F2").Value = "4.98"
$newSheet.Range("G2").Value = "0.5254"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "213002"
$newSheet.Range("C3").Value = "宝盈泛沿海增长混合"
$newSheet.Range("D3").Value = "5.18"
$newSheet.Range("E3").Value = "91.39"
$newSheet.Range("F3").Value = "5.17"
$newSheet.Range("G3").Value = "0.2678"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "000796"
$newSheet.Range("C4").Value = "宝盈睿丰创新灵活配置混合 - C"
$newSheet.Range("D4").Value = "0.64"
$newSheet.Range("E4").Value = "89.70"
$newSheet.Range("F4").Value = "5.56"
$newSheet.Range("G4").Value = "0.0356"
$newSheet.Range("H4").Value = 7

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "000794"
$newSheet.Range("C5").Value = "宝盈睿丰创新灵活配置混合 - A/B"
$newSheet.Range("D5").Value = "0.41"
$newSheet.Range("E5").Value = "89.70"
$newSheet.Range("F5").Value = "5.56"
$newSheet.Range("G5").Value = "0.0228"
$newSheet.Range("H5").Value = 7

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a 2022-Q3 row, pushing the
#    existing 2021-Q4 / 2021-Q2 rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Row 4 is brand new - copy the index-column styling (s=2) used by A2/A3
# onto A4 before writing its value.
$summary.Range("A2").Copy($summary.Range("A4"))

$summary.Range("B4").Value = "2021-Q2"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.01

$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.89

$summary.Range("A4").Value = 2
$summary.Range("A3").Value = 1
$summary.Range("A2").Value = 0

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.85

# ---------------------------------------------------------------------------
# Restore the originally-selected tab ("2021-Q2", the last sheet) as the
# active sheet - Worksheets.Add() above made the freshly-inserted sheet
# active, which would otherwise move the tabSelected flag.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
